$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 34, pushing existing rows 34..87 down to 36..89
$ws.Rows("34:35").Insert()

# Populate new row 34
$ws.Cells.Item(34,1).Value = 1
$ws.Cells.Item(34,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34,3).Value = "Arica y Parinacota"
$ws.Cells.Item(34,4).Value = 44994
$ws.Cells.Item(34,5).Value = 15
$ws.Cells.Item(34,6).Value = "Fruta"
$ws.Cells.Item(34,7).Value = 100103
$ws.Cells.Item(34,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(34,9).Value = 100103004
$ws.Cells.Item(34,10).Value = "Durazno"
$ws.Cells.Item(34,11).Value = "Phillips Cling"
$ws.Cells.Item(34,12).Value = "Segunda"
$ws.Cells.Item(34,13).Value = 550
$ws.Cells.Item(34,14).Value = 20000
$ws.Cells.Item(34,15).Value = 21000
$ws.Cells.Item(34,16).Value = 20545
$ws.Cells.Item(34,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(34,18).Value = "Región de O'Higgins"
$ws.Cells.Item(34,19).Value = 1141
$ws.Cells.Item(34,20).Value = 18

# Populate new row 35
$ws.Cells.Item(35,1).Value = 1
$ws.Cells.Item(35,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(35,3).Value = "Arica y Parinacota"
$ws.Cells.Item(35,4).Value = 44994
$ws.Cells.Item(35,5).Value = 15
$ws.Cells.Item(35,6).Value = "Fruta"
$ws.Cells.Item(35,7).Value = 100103
$ws.Cells.Item(35,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(35,9).Value = 100103004
$ws.Cells.Item(35,10).Value = "Durazno"
$ws.Cells.Item(35,11).Value = "September Sun"
$ws.Cells.Item(35,12).Value = "Segunda"
$ws.Cells.Item(35,13).Value = 250
$ws.Cells.Item(35,14).Value = 20000
$ws.Cells.Item(35,15).Value = 21000
$ws.Cells.Item(35,16).Value = 20600
$ws.Cells.Item(35,17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(35,18).Value = "Región de O'Higgins"
$ws.Cells.Item(35,19).Value = 1144
$ws.Cells.Item(35,20).Value = 18
